# Added handling of common packages.
# Update the "classFields" sheet (field name / field type for the Order
# domain class) so that the generated field metadata matches the actual
# set of fields used once common package handling was added.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("classFields")

# Row layout: A=Class Name, B=Field Name, C=Field Modifier, D=Field Type
# Data rows are 2..8 (row 1 is the header row).

$ws.Range("B2").Value = "productId"
$ws.Range("D2").Value = "java.lang.Long"

$ws.Range("B3").Value = "status"
$ws.Range("D3").Value = "java.lang.String"

$ws.Range("B4").Value = "price"
$ws.Range("D4").Value = "int"

$ws.Range("B5").Value = "productCount"
$ws.Range("D5").Value = "int"

$ws.Range("B6").Value = "source"
$ws.Range("D6").Value = "java.lang.String"

$ws.Range("B7").Value = "customerId"
$ws.Range("D7").Value = "java.lang.Long"

$ws.Range("B8").Value = "id"
$ws.Range("D8").Value = "java.lang.Long"
